$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the misspelled director name and drop the extra rows; only the
# "Director" header (A1) and the single corrected entry (A2) should remain.
$ws.Range("A2").Value = "Steven Spielberg"
$ws.Range("A3:A6").ClearContents()

[void]$ws.Range("A2").Select()
